$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3432
$ws.Range("I76").Value = 2500
$ws.Range("J76").Value = 3548.5
$ws.Range("K76").Value = 2500
$ws.Range("L76").Value = 3548.5
$ws.Range("M76").Value = -2185
$ws.Range("N76").Value = -4178.5
$ws.Range("H79").Value = 3432
$ws.Range("I79").Value = 2500
$ws.Range("J79").Value = 3548.5
$ws.Range("K79").Value = 2500
$ws.Range("L79").Value = 3548.5
$ws.Range("M79").Value = -1408
$ws.Range("N79").Value = -5732.5
$ws.Range("H86").Value = 1734.2106
$ws.Range("I86").Value = 1675
$ws.Range("J86").Value = 1835.7142
$ws.Range("K86").Value = 1675
$ws.Range("L86").Value = 1835.7142
$ws.Range("M86").Value = -552
$ws.Range("N86").Value = -4081.7142
$ws.Range("H89").Value = 1734.2106
$ws.Range("I89").Value = 1675
$ws.Range("J89").Value = 1835.7142
$ws.Range("K89").Value = 8375
$ws.Range("L89").Value = 9178.571
$ws.Range("M89").Value = -2759
$ws.Range("N89").Value = -20410.571
$ws.Range("H132").Value = 2122.2793
$ws.Range("I132").Value = 1741.362
$ws.Range("J132").Value = 4331.6
$ws.Range("K132").Value = 5224.086
$ws.Range("L132").Value = 12994.8
$ws.Range("M132").Value = -2694.086
$ws.Range("N132").Value = -18054.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 8000
$ws.Range("I31").Value = 8000
$ws.Range("K31").Value = 8000
$ws.Range("M31").Value = -7706
$ws.Range("H32").Value = 10320291
$ws.Range("I32").Value = 12990470
$ws.Range("K32").Value = 12990470
$ws.Range("M32").Value = -12990183
$ws.Range("H97").Value = 912.3333
$ws.Range("I97").Value = 640.4167
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 640.4167
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -144.4167
$ws.Range("N97").Value = -2992
$ws.Range("H128").Value = 32500
$ws.Range("J128").Value = 32500
$ws.Range("L128").Value = 32500
$ws.Range("N128").Value = -42460
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2323
$ws.Range("I86").Value = 2039.8
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 2039.8
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -916.8
$ws.Range("N86").Value = -4746
$ws.Range("H89").Value = 2323
$ws.Range("I89").Value = 2039.8
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 10199
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -4583
$ws.Range("N89").Value = -23732
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 21666.666
$ws.Range("J82").Value = 21666.666
$ws.Range("L82").Value = 21666.666
$ws.Range("N82").Value = -22388.666
$ws.Range("H85").Value = 21666.666
$ws.Range("J85").Value = 21666.666
$ws.Range("L85").Value = 21666.666
$ws.Range("N85").Value = -24162.666
$ws.Range("H102").Value = 25000
$ws.Range("J102").Value = 25000
$ws.Range("L102").Value = 25000
$ws.Range("N102").Value = -29868
$ws.Range("H105").Value = 544.1429000000001
$ws.Range("I105").Value = 544.1429000000001
$ws.Range("K105").Value = 544.1429000000001
$ws.Range("M105").Value = 1202.8571
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 99.833336
$ws.Range("I8").Value = 99.833336
$ws.Range("K8").Value = 299.500008
$ws.Range("M8").Value = -160.500008
$ws.Range("H132").Value = 9304929
$ws.Range("I132").Value = 397.75
$ws.Range("J132").Value = 11786137
$ws.Range("K132").Value = 3579.75
$ws.Range("L132").Value = 106075233
$ws.Range("M132").Value = -1049.75
$ws.Range("N132").Value = -106080293
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 14283.333
$ws.Range("J15").Value = 14283.333
$ws.Range("L15").Value = 14283.333
$ws.Range("N15").Value = -14859.333
$ws.Range("H80").Value = 2983.7058
$ws.Range("I80").Value = 2764.6428
$ws.Range("K80").Value = 2764.6428
$ws.Range("M80").Value = -1766.6428
$ws.Range("H81").Value = 14283.333
$ws.Range("J81").Value = 14283.333
$ws.Range("L81").Value = 14283.333
$ws.Range("N81").Value = -16279.333
$ws.Range("H83").Value = 2983.7058
$ws.Range("I83").Value = 2764.6428
$ws.Range("K83").Value = 13823.214
$ws.Range("M83").Value = -8831.214
$ws.Range("H84").Value = 14283.333
$ws.Range("J84").Value = 14283.333
$ws.Range("L84").Value = 42849.999
$ws.Range("N84").Value = -52833.999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3065.875
$ws.Range("I82").Value = 2589.3333
$ws.Range("J82").Value = 3678.5715
$ws.Range("K82").Value = 2589.3333
$ws.Range("L82").Value = 3678.5715
$ws.Range("M82").Value = -2228.3333
$ws.Range("N82").Value = -4400.5715
$ws.Range("H85").Value = 3065.875
$ws.Range("I85").Value = 2589.3333
$ws.Range("J85").Value = 3678.5715
$ws.Range("K85").Value = 2589.3333
$ws.Range("L85").Value = 3678.5715
$ws.Range("M85").Value = -1341.3333
$ws.Range("N85").Value = -6174.5715
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11124083
$ws.Range("I62").Value = 50050000
$ws.Range("J62").Value = 2392.8572
$ws.Range("K62").Value = 50050000
$ws.Range("L62").Value = 2392.8572
$ws.Range("M62").Value = -50049376
$ws.Range("N62").Value = -3640.8572
$ws.Range("H65").Value = 11124083
$ws.Range("I65").Value = 50050000
$ws.Range("J65").Value = 2392.8572
$ws.Range("K65").Value = 250250000
$ws.Range("L65").Value = 11964.286
$ws.Range("M65").Value = -250246880
$ws.Range("N65").Value = -18204.286
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 10000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 10000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -9685
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 10000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 10000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8908
$ws.Range("N73").ClearContents()